$d = $word.ActiveDocument

# 1. Remove "representative_position} {" from the first occurrence
#    (the table cell that originally reads:
#     "{representative_position} {representative_name} {representative_code}")
#    leaving "{representative_name} {representative_code}"
$d.Content.Find.Execute("{representative_position} {representative_name}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{representative_name}", 2)

# 2. Mark the "Numatytasispastraiposriftas" (Default Paragraph Font) style as
#    semi-hidden.
$style = $d.Styles("Numatytasispastraiposriftas")
$style.Hidden = $true
